$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# Row 30 is currently the "freshly added" event row: A30/L30:O30 carry the
# event data, while B30:K30 and P30 are blank placeholder cells.
#
# This edit (per the commit) adds ANOTHER new event row (31) with the same
# shape row 30 used to have (blank B:K / P, same date/event/correction/by
# values), and normalizes row 30's old blanks to the literal text "nan" -
# matching every earlier row (2-29) on this sheet, which already store
# their blanks as the text "nan" rather than as empty cells.

# 1) Clone row 30 (as it stands right now) down into row 31. This keeps
#    A31 typed as text (matching A30) and leaves B31:K31/P31 blank, exactly
#    like row 30 was before this edit.
$ws.Range("A30:P30").Copy($ws.Range("A31:P31"))

# 2) Normalize row 30's blank cells to literal "nan" text.
$naCols = @("B","C","D","E","F","G","H","I","J","K","P")
foreach ($col in $naCols) {
    $ws.Range($col + "30").Value = "nan"
}
